# Weekly data refresh: a new daily price record for
# "Femacal de La Calera - Acelga" is inserted as row 534 (the sheet is kept
# in reverse-chronological-ish order by the source feed), pushing the
# existing rows 534-652 down to 535-653.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 534; Excel shifts rows 534..652 down
# to 535..653 and the new row inherits formatting from the row above it
# (keeps the date-style on column D, etc.)
$ws.Rows(534).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(534, 1).Value  = 3
$ws.Cells.Item(534, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(534, 3).Value  = "Coquimbo"
$ws.Cells.Item(534, 4).Value  = 45244
$ws.Cells.Item(534, 5).Value  = 5
$ws.Cells.Item(534, 6).Value  = 100112009
$ws.Cells.Item(534, 7).Value  = "Acelga"
$ws.Cells.Item(534, 8).Value  = "Sin especificar"
$ws.Cells.Item(534, 9).Value  = "Primera"
$ws.Cells.Item(534, 10).Value = 210
$ws.Cells.Item(534, 11).Value = 3300
$ws.Cells.Item(534, 12).Value = 3500
$ws.Cells.Item(534, 13).Value = 3414
$ws.Cells.Item(534, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(534, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(534, 16).Value = 569
$ws.Cells.Item(534, 17).Value = 6
$ws.Cells.Item(534, 18).Value = "Hortaliza"
